# Alpha3F-HW30 notebook rerun: two new HKL rows ("Holden", "Rizzie Spiral")
# inserted after "Spiral5", and "Thomas Hex" renamed to "Matthies Hex".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new blank rows at 4:5 (pushes old rows 4-29 down to 6-31) ---
$ws.Rows("4:5").Insert()

# Copy the formatting (style) of row 3 (A/B columns) onto the two new rows so
# the new "index" column (A) keeps the same bold/bordered/centered style as
# every other data row, instead of the auto-derived style Insert() produces.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B5").PasteSpecial(-4122)

# --- Fill in the two new data rows ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 0.9844814532631846
$ws.Range("D4").Value = 1.021681091934135
$ws.Range("E4").Value = 1.004447507156223
$ws.Range("F4").Value = 1.006948378287868
$ws.Range("G4").Value = 0.9730003660953098
$ws.Range("H4").Value = 1.007977542516593
$ws.Range("I4").Value = 0.9730003660953098
$ws.Range("J4").Value = 1.004447507156223
$ws.Range("K4").Value = 1.004447507156223
$ws.Range("L4").Value = 1.007977542516593
$ws.Range("M4").Value = 0.9904889543059516
$ws.Range("N4").Value = 0.9904889543059516
$ws.Range("O4").Value = 0.9884864539583625
$ws.Range("P4").Value = 0.995141805256042
$ws.Range("Q4").Value = 0.995141805256042
$ws.Range("R4").Value = 0.9974682307310871
$ws.Range("S4").Value = 0.9974682307310871
$ws.Range("T4").Value = 0.999756056542219

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 0.8880291884488235
$ws.Range("D5").Value = 1.104816805048931
$ws.Range("E5").Value = 1.078078843091477
$ws.Range("F5").Value = 1.052230291464663
$ws.Range("G5").Value = 0.7819310930928848
$ws.Range("H5").Value = 1.041593057705963
$ws.Range("I5").Value = 0.7819310930928848
$ws.Range("J5").Value = 1.078078843091477
$ws.Range("K5").Value = 1.078078843091477
$ws.Range("L5").Value = 1.041593057705963
$ws.Range("M5").Value = 0.9117620753994238
$ws.Range("N5").Value = 0.9117620753994238
$ws.Range("O5").Value = 0.9038511130825571
$ws.Range("P5").Value = 0.9672009979634414
$ws.Range("Q5").Value = 0.9672009979634416
$ws.Range("R5").Value = 0.9949204592454504
$ws.Range("S5").Value = 0.9949204592454504
$ws.Range("T5").Value = 0.9911132131421238

# --- Rename "Thomas Hex" -> "Matthies Hex" (now on row 11 after the shift) ---
$ws.Range("B11").Value = "Matthies Hex"
